$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row40 (hunk 0)
$ws_ALC.Range("H40").Value = 3211.111
$ws_ALC.Range("I40").Value = 2225
$ws_ALC.Range("J40").Value = 4000
$ws_ALC.Range("K40").Value = 2225
$ws_ALC.Range("L40").Value = 4000
$ws_ALC.Range("M40").Value = -2050
$ws_ALC.Range("N40").Value = -4350

# ALC!row62 (hunk 1)
$ws_ALC.Range("H62").Value = 76493.71000000001
$ws_ALC.Range("I62").Value = 146773.28
$ws_ALC.Range("J62").Value = 6214.143
$ws_ALC.Range("K62").Value = 146773.28
$ws_ALC.Range("L62").Value = 6214.143
$ws_ALC.Range("M62").Value = -146149.28
$ws_ALC.Range("N62").Value = -7462.143

# ALC!row65 (hunk 2)
$ws_ALC.Range("H65").Value = 76493.71000000001
$ws_ALC.Range("I65").Value = 146773.28
$ws_ALC.Range("J65").Value = 6214.143
$ws_ALC.Range("K65").Value = 733866.4
$ws_ALC.Range("L65").Value = 31070.715
$ws_ALC.Range("M65").Value = -730746.4
$ws_ALC.Range("N65").Value = -37310.715

# ALC!row69 (hunk 3)
$ws_ALC.Range("H69").Value = 0
$ws_ALC.Range("J69").Value = 0
$ws_ALC.Range("L69").Value = 0
$ws_ALC.Range("N69").ClearContents()

# ALC!row72 (hunk 4)
$ws_ALC.Range("H72").Value = 0
$ws_ALC.Range("J72").Value = 0
$ws_ALC.Range("L72").Value = 0
$ws_ALC.Range("N72").ClearContents()

# ALC!row76 (hunk 5)
$ws_ALC.Range("H76").Value = 3654
$ws_ALC.Range("I76").Value = 3317.5
$ws_ALC.Range("K76").Value = 3317.5
$ws_ALC.Range("M76").Value = -3002.5

# ALC!row79 (hunk 6)
$ws_ALC.Range("H79").Value = 3654
$ws_ALC.Range("I79").Value = 3317.5
$ws_ALC.Range("K79").Value = 3317.5
$ws_ALC.Range("M79").Value = -2225.5

# ALC!row103 (hunk 7)
$ws_ALC.Range("H103").Value = 367.30768
$ws_ALC.Range("I103").Value = 240
$ws_ALC.Range("J103").Value = 540.9091
$ws_ALC.Range("K103").Value = 720
$ws_ALC.Range("L103").Value = 1622.7273
$ws_ALC.Range("M103").Value = -134
$ws_ALC.Range("N103").Value = -2794.7273

# ALC!row104 (hunk 8)
$ws_ALC.Range("H104").Value = 1162.25
$ws_ALC.Range("I104").Value = 883.3333
$ws_ALC.Range("K104").Value = 2649.9999
$ws_ALC.Range("M104").Value = -902.9998999999998

# ALC!row107 (hunk 9)
$ws_ALC.Range("H107").Value = 904.8823
$ws_ALC.Range("I107").Value = 908.4375
$ws_ALC.Range("K107").Value = 908.4375
$ws_ALC.Range("M107").Value = 1011.5625

# ALC!row113 (hunk 10)
$ws_ALC.Range("H113").Value = 5240
$ws_ALC.Range("J113").Value = 6000
$ws_ALC.Range("L113").Value = 6000
$ws_ALC.Range("N113").Value = -12508

# ALC!row117 (hunk 11)
$ws_ALC.Range("H117").Value = 0
$ws_ALC.Range("I117").Value = 0
$ws_ALC.Range("K117").Value = 0
$ws_ALC.Range("M117").ClearContents()

# ARM!row38 (hunk 12)
$ws_ARM.Range("H38").Value = 0
$ws_ARM.Range("I38").Value = 0
$ws_ARM.Range("J38").Value = 0
$ws_ARM.Range("K38").Value = 0
$ws_ARM.Range("L38").Value = 0
$ws_ARM.Range("M38").ClearContents()
$ws_ARM.Range("N38").ClearContents()

# ARM!row63 (hunk 13)
$ws_ARM.Range("H63").Value = 2634.4
$ws_ARM.Range("I63").Value = 2371.5557
$ws_ARM.Range("K63").Value = 2371.5557
$ws_ARM.Range("M63").Value = -1685.5557

# ARM!row66 (hunk 14)
$ws_ARM.Range("H66").Value = 2634.4
$ws_ARM.Range("I66").Value = 2371.5557
$ws_ARM.Range("K66").Value = 11857.7785
$ws_ARM.Range("M66").Value = -8425.7785

# ARM!row74 (hunk 15)
$ws_ARM.Range("H74").Value = 1983.5
$ws_ARM.Range("I74").Value = 2007.1
$ws_ARM.Range("K74").Value = 2007.1
$ws_ARM.Range("M74").Value = -1133.1

# ARM!row77 (hunk 16)
$ws_ARM.Range("H77").Value = 1983.5
$ws_ARM.Range("I77").Value = 2007.1
$ws_ARM.Range("K77").Value = 10035.5
$ws_ARM.Range("M77").Value = -5667.5

# ARM!row97 (hunk 17)
$ws_ARM.Range("H97").Value = 2723
$ws_ARM.Range("I97").Value = 1592.5834
$ws_ARM.Range("J97").Value = 9505.5
$ws_ARM.Range("K97").Value = 1592.5834
$ws_ARM.Range("L97").Value = 9505.5
$ws_ARM.Range("M97").Value = -1096.5834
$ws_ARM.Range("N97").Value = -10497.5

# BSM!row82 (hunk 18)
$ws_BSM.Range("H82").Value = 20919
$ws_BSM.Range("J82").Value = 25000
$ws_BSM.Range("L82").Value = 25000
$ws_BSM.Range("N82").Value = -25766

# BSM!row85 (hunk 19)
$ws_BSM.Range("H85").Value = 20919
$ws_BSM.Range("J85").Value = 25000
$ws_BSM.Range("L85").Value = 25000
$ws_BSM.Range("N85").Value = -27652

# BSM!row94 (hunk 20)
$ws_BSM.Range("H94").Value = 1121.8182
$ws_BSM.Range("I94").Value = 1130.125
$ws_BSM.Range("J94").Value = 1099.6666
$ws_BSM.Range("K94").Value = 1130.125
$ws_BSM.Range("L94").Value = 1099.6666
$ws_BSM.Range("M94").Value = -679.125
$ws_BSM.Range("N94").Value = -2001.6666

# BSM!row99 (hunk 21)
$ws_BSM.Range("H99").Value = 1505.3636
$ws_BSM.Range("I99").Value = 1510.2
$ws_BSM.Range("J99").Value = 1457
$ws_BSM.Range("K99").Value = 1510.2
$ws_BSM.Range("L99").Value = 1457
$ws_BSM.Range("M99").Value = -12.20000000000005
$ws_BSM.Range("N99").Value = -4453

# BSM!row134 (hunk 22)
$ws_BSM.Range("H134").Value = 1484
$ws_BSM.Range("I134").Value = 1494.2222
$ws_BSM.Range("K134").Value = 4482.6666
$ws_BSM.Range("M134").Value = -1947.6666

# CRP!row31 (hunk 23)
$ws_CRP.Range("H31").Value = 2627.5938
$ws_CRP.Range("I31").Value = 2020.4375
$ws_CRP.Range("J31").Value = 3234.75
$ws_CRP.Range("K31").Value = 2020.4375
$ws_CRP.Range("L31").Value = 3234.75
$ws_CRP.Range("M31").Value = -1725.4375
$ws_CRP.Range("N31").Value = -3824.75

# CRP!row34 (hunk 24)
$ws_CRP.Range("H34").Value = 2627.5938
$ws_CRP.Range("I34").Value = 2020.4375
$ws_CRP.Range("J34").Value = 3234.75
$ws_CRP.Range("K34").Value = 2020.4375
$ws_CRP.Range("L34").Value = 3234.75
$ws_CRP.Range("M34").Value = -1818.4375
$ws_CRP.Range("N34").Value = -3638.75

# CRP!row86 (hunk 25)
$ws_CRP.Range("H86").Value = 5322.7144
$ws_CRP.Range("I86").Value = 4316
$ws_CRP.Range("J86").Value = 6665
$ws_CRP.Range("K86").Value = 4316
$ws_CRP.Range("L86").Value = 6665
$ws_CRP.Range("M86").Value = -3193
$ws_CRP.Range("N86").Value = -8911

# CRP!row89 (hunk 26)
$ws_CRP.Range("H89").Value = 5322.7144
$ws_CRP.Range("I89").Value = 4316
$ws_CRP.Range("J89").Value = 6665
$ws_CRP.Range("K89").Value = 21580
$ws_CRP.Range("L89").Value = 33325
$ws_CRP.Range("M89").Value = -15964
$ws_CRP.Range("N89").Value = -44557

# CRP!row99 (hunk 27)
$ws_CRP.Range("H99").Value = 1999.3334
$ws_CRP.Range("I99").Value = 1999.3334
$ws_CRP.Range("J99").Value = 0
$ws_CRP.Range("K99").Value = 1999.3334
$ws_CRP.Range("L99").Value = 0
$ws_CRP.Range("M99").Value = -501.3334
$ws_CRP.Range("N99").ClearContents()

# CRP!row105 (hunk 28)
$ws_CRP.Range("H105").Value = 8052.933
$ws_CRP.Range("I105").Value = 6866
$ws_CRP.Range("K105").Value = 6866
$ws_CRP.Range("M105").Value = -5119

# CRP!row126 (hunk 29)
$ws_CRP.Range("H126").Value = 1999.3334
$ws_CRP.Range("I126").Value = 1999.3334
$ws_CRP.Range("J126").Value = 0
$ws_CRP.Range("K126").Value = 5998.0002
$ws_CRP.Range("L126").Value = 0
$ws_CRP.Range("M126").Value = -3528.0002
$ws_CRP.Range("N126").ClearContents()

# CRP!row132 (hunk 30)
$ws_CRP.Range("H132").Value = 2080.7856
$ws_CRP.Range("I132").Value = 2089.95
$ws_CRP.Range("J132").Value = 1897.5
$ws_CRP.Range("K132").Value = 6269.849999999999
$ws_CRP.Range("L132").Value = 5692.5
$ws_CRP.Range("M132").Value = -3739.849999999999
$ws_CRP.Range("N132").Value = -10752.5

# CRP!row134 (hunk 31)
$ws_CRP.Range("H134").Value = 2829.027
$ws_CRP.Range("I134").Value = 2649.9355
$ws_CRP.Range("J134").Value = 3754.3333
$ws_CRP.Range("K134").Value = 7949.806500000001
$ws_CRP.Range("L134").Value = 11262.9999
$ws_CRP.Range("M134").Value = -5414.806500000001
$ws_CRP.Range("N134").Value = -16332.9999

# CUL!row23 (hunk 32)
$ws_CUL.Range("H23").Value = 98.333336
$ws_CUL.Range("I23").Value = 95
$ws_CUL.Range("K23").Value = 285
$ws_CUL.Range("M23").Value = -50

# CUL!row64 (hunk 33)
$ws_CUL.Range("H64").Value = 5217.375
$ws_CUL.Range("J64").Value = 5217.375
$ws_CUL.Range("L64").Value = 15652.125
$ws_CUL.Range("N64").Value = -16192.125

# CUL!row67 (hunk 34)
$ws_CUL.Range("H67").Value = 5217.375
$ws_CUL.Range("J67").Value = 5217.375
$ws_CUL.Range("L67").Value = 15652.125
$ws_CUL.Range("N67").Value = -17524.125

# CUL!row115 (hunk 35)
$ws_CUL.Range("H115").Value = 1604.6666
$ws_CUL.Range("I115").Value = 907
$ws_CUL.Range("K115").Value = 2721
$ws_CUL.Range("M115").Value = -1546

# CUL!row120 (hunk 36)
$ws_CUL.Range("H120").Value = 19722.8
$ws_CUL.Range("J120").Value = 20000
$ws_CUL.Range("L120").Value = 60000
$ws_CUL.Range("N120").Value = -69676

# CUL!row128 (hunk 37)
$ws_CUL.Range("H128").Value = 246245.14
$ws_CUL.Range("I128").Value = 246245.14
$ws_CUL.Range("K128").Value = 738735.42
$ws_CUL.Range("M128").Value = -733755.42

# CUL!row131 (hunk 38)
$ws_CUL.Range("H131").Value = 5936.28
$ws_CUL.Range("I131").Value = 13864.375
$ws_CUL.Range("J131").Value = 2205.4119
$ws_CUL.Range("K131").Value = 41593.125
$ws_CUL.Range("L131").Value = 6616.2357
$ws_CUL.Range("M131").Value = -36553.125
$ws_CUL.Range("N131").Value = -16696.2357

# GSM!row102 (hunk 39)
$ws_GSM.Range("H102").Value = 2166.2307
$ws_GSM.Range("I102").Value = 1806.2
$ws_GSM.Range("K102").Value = 1806.2
$ws_GSM.Range("M102").Value = -184.2

# LTW!row7 (hunk 40)
$ws_LTW.Range("H7").Value = 8284.125
$ws_LTW.Range("I7").Value = 7295.9165
$ws_LTW.Range("K7").Value = 7295.9165
$ws_LTW.Range("M7").Value = -7183.9165

# LTW!row16 (hunk 41)
$ws_LTW.Range("H16").Value = 1492.1818
$ws_LTW.Range("I16").Value = 1601.4
$ws_LTW.Range("J16").Value = 400
$ws_LTW.Range("K16").Value = 1601.4
$ws_LTW.Range("L16").Value = 400
$ws_LTW.Range("M16").Value = -1431.4
$ws_LTW.Range("N16").Value = -740

# LTW!row40 (hunk 42)
$ws_LTW.Range("H40").Value = 5401
$ws_LTW.Range("I40").Value = 5098.143
$ws_LTW.Range("K40").Value = 5098.143
$ws_LTW.Range("M40").Value = -4962.143

# LTW!row100 (hunk 43)
$ws_LTW.Range("H100").Value = 24676.688
$ws_LTW.Range("I100").Value = 3939.8
$ws_LTW.Range("K100").Value = 3939.8
$ws_LTW.Range("M100").Value = -3398.8

# LTW!row107 (hunk 44)
$ws_LTW.Range("H107").Value = 2213.4285
$ws_LTW.Range("I107").Value = 2213.4285
$ws_LTW.Range("K107").Value = 2213.4285
$ws_LTW.Range("M107").Value = -293.4285

# LTW!row126 (hunk 45)
$ws_LTW.Range("H126").Value = 8284.125
$ws_LTW.Range("I126").Value = 7295.9165
$ws_LTW.Range("K126").Value = 21887.7495
$ws_LTW.Range("M126").Value = -19417.7495

# WVR!row98 (hunk 46)
$ws_WVR.Range("H98").Value = 30295
$ws_WVR.Range("J98").Value = 30295
$ws_WVR.Range("L98").Value = 30295
$ws_WVR.Range("N98").Value = -36285
